$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.657.71'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '3.650.67'
$ws.Range('E3').Value = '  +1.80%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '243.03'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').Value = '1.81'
$ws.Range('E6').Value = '  +16.29%  '
$ws.Range('D7').Value = '656.52'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '0.423'
$ws.Range('E8').Value = '  +4.04%  '
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('D10').Value = '0.999'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').Value = '3.648.89'
$ws.Range('E11').Value = '  +1.84%  '
$ws.Range('D12').Value = '44.45'
$ws.Range('E12').Value = '  +2.70%  '
$ws.Range('D13').Value = '0.204'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').Value = '4.326.66'
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').Value = '96.408.99'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '0.0000259'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '8.15'
$ws.Range('E18').Value = '  +5.14%  '
$ws.Range('D19').Value = '3.647.03'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').Value = '13.20'
$ws.Range('E20').Value = '  +4.75%  '
$ws.Range('D21').Value = '18.50'
$ws.Range('E21').Value = '  +3.88%  '
$ws.Range('D22').Value = '0.529'
$ws.Range('E22').Value = '  +7.68%  '
$ws.Range('D23').Value = '3.45'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '512.81'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = '0.0000205'
$ws.Range('E25').Value = '  +2.13%  '
$ws.Range('D26').Value = '6.90'
$ws.Range('E26').Value = '  +1.15%  '
$ws.Range('D27').Value = '100.88'
$ws.Range('E27').Value = '  +4.32%  '
$ws.Range('D28').Value = '13.21'
$ws.Range('E28').Value = '  +3.59%  '
$ws.Range('D29').Value = '0.165'
$ws.Range('E29').Value = '  +10.85%  '
$ws.Range('D30').Value = '3.05'
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('D31').Value = '11.94'
$ws.Range('E31').Value = '  +3.91%  '
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').Value = '0.185'
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('D34').Value = '33.39'
$ws.Range('E34').Value = '  +5.50%  '
$ws.Range('E35').Value = '  +0.53%  '
$ws.Range('D36').Value = '1.72'
$ws.Range('E36').Value = '  +7.96%  '
$ws.Range('D37').Value = '0.581'
$ws.Range('E37').Value = '  +2.85%  '
$ws.Range('D38').Value = '8.88'
$ws.Range('E38').Value = '  +4.61%  '
$ws.Range('D39').Value = '617.00'
$ws.Range('E39').Value = '  +3.06%  '
$ws.Range('D40').Value = '43.10'
$ws.Range('E40').Value = '  +25.44%  '
$ws.Range('E41').Value = '  +3.16%  '
$ws.Range('D42').Value = '0.959'
$ws.Range('E42').Value = '  +5.61%  '
$ws.Range('D43').Value = '1.95'
$ws.Range('E43').Value = '  +6.03%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '6.15'
$ws.Range('E45').Value = '  +7.17%  '
$ws.Range('D46').Value = '0.0441'
$ws.Range('E46').Value = '  +4.89%  '
$ws.Range('D47').Value = '2.31'
$ws.Range('E47').Value = '  +1.42%  '
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('D49').Value = '8.55'
$ws.Range('E49').Value = '  +3.87%  '
$ws.Range('D50').Value = '0.402'
$ws.Range('E50').Value = '  +14.37%  '
$ws.Range('D51').Value = '54.52'
$ws.Range('E51').Value = '  +1.72%  '
